# Auto-generated script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.553.45"
$ws.Range("E2").Value = "  -1.84%  "
$ws.Range("D3").Value = "2.971.78"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "539.16"
$ws.Range("E5").Value = "  -2.31%  "
$ws.Range("D6").Value = "149.24"
$ws.Range("E6").Value = "  -3.27%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "0.562"
$ws.Range("E8").Value = "  +2.27%  "
$ws.Range("D9").Value = "2.982.97"
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("D10").Value = "0.112"
$ws.Range("E10").Value = "  +1.06%  "
$ws.Range("E11").Value = "  -4.49%  "
$ws.Range("D12").Value = "0.365"
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("D13").Value = "3.496.15"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("E14").Value = "  +1.68%  "
$ws.Range("D15").Value = "61.613.83"
$ws.Range("E15").Value = "  -1.93%  "
$ws.Range("D16").Value = "23.75"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "2.982.53"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("E18").Value = "  -1.62%  "
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("D20").Value = "11.94"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "373.38"
$ws.Range("E21").Value = "  -3.76%  "
$ws.Range("E22").Value = "  +1.23%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "65.61"
$ws.Range("E24").Value = "  +1.19%  "
$ws.Range("D25").Value = "3.109.86"
$ws.Range("E25").Value = "  -1.62%  "
$ws.Range("D26").Value = "0.466"
$ws.Range("E26").Value = "  +1.36%  "
$ws.Range("D27").Value = "0.188"
$ws.Range("E27").Value = "  +2.17%  "
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("D29").Value = "0.0₃0907"
$ws.Range("E29").Value = "  -5.52%  "
$ws.Range("D30").Value = "8.13"
$ws.Range("E30").Value = "  -5.28%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").Value = "1.71"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("D34").Value = "159.91"
$ws.Range("E34").Value = "  -0.81%  "
$ws.Range("E35").Value = "  -2.50%  "
$ws.Range("E36").Value = "  -1.72%  "
$ws.Range("E37").Value = "  -3.82%  "
$ws.Range("D38").Value = "1.25"
$ws.Range("E38").Value = "  -2.57%  "
$ws.Range("E39").Value = "  -3.26%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "37.21"
$ws.Range("E40").Value = "  -0.66%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.399.33"
$ws.Range("E41").Value = "  -4.94%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "0.668"
$ws.Range("E43").Value = "  +1.49%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "21.79"
$ws.Range("E44").Value = "  -2.78%  "
$ws.Range("D45").Value = "0.0585"
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("E46").Value = "  +4.17%  "
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("D49").Value = "267.30"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").Value = "0.0944"
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("D51").Value = "10.40"
$ws.Range("E51").Value = "  -0.91%  "
